$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a new table row above the physical row 55 (ListRows.Item(47)),
# shifting existing rows 55-133 down to 56-134, and grow the table to A8:K134.
$ws.Rows.Item(55).Insert()
$tbl.Resize($ws.Range("A8:K134"))

# The freshly inserted row inherits blank/minimal formatting; copy the
# normal data-row formatting (borders, number formats, etc.) from the row
# directly below (now row 56, a normal interior table row) onto row 55.
$ws.Range("A56:K56").Copy()
$ws.Range("A55:K55").PasteSpecial(-4122)

# Re-assert the calculated-column formulas for the row that was inserted
# and for the row that got pushed down to the bottom of the table so they
# evaluate cleanly instead of carrying over a stale/broken reference.
$ws.Range("G55").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G134").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Fill in the new leave-card entry on row 55.
$ws.Range("B55").Value = "SL(3-0-0)"
$ws.Range("H55").Value = 3
$ws.Range("K55").Value = "5/8-10/2023"

# Move the active selection to where it ended up after the insert (one row down).
$ws.Range("B56").Select()
